$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1165.3846
$ws.Range("I55").Value = 1416.4445
$ws.Range("J55").Value = 600.5
$ws.Range("K55").Value = 1416.4445
$ws.Range("L55").Value = 600.5
$ws.Range("M55").Value = -1202.4445
$ws.Range("N55").Value = -1028.5

$ws.Range("H82").Value = 2781.5386
$ws.Range("I82").Value = 522.8570999999999
$ws.Range("J82").Value = 5416.6665
$ws.Range("K82").Value = 1568.5713
$ws.Range("L82").Value = 16249.9995
$ws.Range("M82").Value = -1162.5713
$ws.Range("N82").Value = -17061.9995

$ws.Range("H85").Value = 2781.5386
$ws.Range("I85").Value = 522.8570999999999
$ws.Range("J85").Value = 5416.6665
$ws.Range("K85").Value = 1568.5713
$ws.Range("L85").Value = 16249.9995
$ws.Range("M85").Value = -164.5712999999998
$ws.Range("N85").Value = -19057.9995

$ws.Range("H135").Value = 774.3461
$ws.Range("I135").Value = 769.43475
$ws.Range("K135").Value = 6924.91275
$ws.Range("M135").Value = -4389.91275

$ws.Range("H137").Value = 2527.4707
$ws.Range("I137").Value = 2801.1072
$ws.Range("J137").Value = 2194.348
$ws.Range("K137").Value = 8403.321599999999
$ws.Range("L137").Value = 6583.044
$ws.Range("M137").Value = -5853.321599999999
$ws.Range("N137").Value = -11683.044

$ws.Range("H138").Value = 3233.698
$ws.Range("I138").Value = 1474.4546
$ws.Range("J138").Value = 6136.45
$ws.Range("K138").Value = 4423.3638
$ws.Range("L138").Value = 18409.35
$ws.Range("M138").Value = 716.6361999999999
$ws.Range("N138").Value = -28689.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1402.8368
$ws.Range("I45").Value = 942.1389
$ws.Range("J45").Value = 2678.6155
$ws.Range("K45").Value = 942.1389
$ws.Range("L45").Value = 2678.6155
$ws.Range("M45").Value = -565.1389
$ws.Range("N45").Value = -3432.6155

$ws.Range("H61").Value = 3294.4167
$ws.Range("I61").Value = 1228.7273
$ws.Range("J61").Value = 5042.3076
$ws.Range("K61").Value = 1228.7273
$ws.Range("L61").Value = 5042.3076
$ws.Range("M61").Value = -1016.7273
$ws.Range("N61").Value = -5466.3076

$ws.Range("H74").Value = 648.7692
$ws.Range("I74").Value = 554.9
$ws.Range("J74").Value = 961.6667
$ws.Range("K74").Value = 554.9
$ws.Range("L74").Value = 961.6667
$ws.Range("M74").Value = 319.1
$ws.Range("N74").Value = -2709.6667

$ws.Range("H77").Value = 648.7692
$ws.Range("I77").Value = 554.9
$ws.Range("J77").Value = 961.6667
$ws.Range("K77").Value = 2774.5
$ws.Range("L77").Value = 4808.3335
$ws.Range("M77").Value = 1593.5
$ws.Range("N77").Value = -13544.3335

$ws.Range("H102").Value = 2693.8125
$ws.Range("I102").Value = 2212.1
$ws.Range("J102").Value = 3496.6667
$ws.Range("K102").Value = 2212.1
$ws.Range("L102").Value = 3496.6667
$ws.Range("M102").Value = -590.0999999999999
$ws.Range("N102").Value = -6740.6667

$ws.Range("H132").Value = 29415530
$ws.Range("I132").Value = 41670340
$ws.Range("K132").Value = 125011020
$ws.Range("M132").Value = -125008490

$ws.Range("H136").Value = 3294.4167
$ws.Range("I136").Value = 1228.7273
$ws.Range("J136").Value = 5042.3076
$ws.Range("K136").Value = 3686.1819
$ws.Range("L136").Value = 15126.9228
$ws.Range("M136").Value = -1136.1819
$ws.Range("N136").Value = -20226.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4170.6
$ws.Range("I86").Value = 4081.2
$ws.Range("J86").Value = 4260
$ws.Range("K86").Value = 4081.2
$ws.Range("L86").Value = 4260
$ws.Range("M86").Value = -2958.2
$ws.Range("N86").Value = -6506

$ws.Range("H89").Value = 4170.6
$ws.Range("I89").Value = 4081.2
$ws.Range("J89").Value = 4260
$ws.Range("K89").Value = 20406
$ws.Range("L89").Value = 21300
$ws.Range("M89").Value = -14790
$ws.Range("N89").Value = -32532

$ws.Range("H134").Value = 2525.8572
$ws.Range("I134").Value = 1360.1818
$ws.Range("J134").Value = 6800
$ws.Range("K134").Value = 4080.5454
$ws.Range("L134").Value = 20400
$ws.Range("M134").Value = -1545.5454
$ws.Range("N134").Value = -25470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2273.1462
$ws.Range("I31").Value = 1466.44
$ws.Range("J31").Value = 3533.625
$ws.Range("K31").Value = 1466.44
$ws.Range("L31").Value = 3533.625
$ws.Range("M31").Value = -1171.44
$ws.Range("N31").Value = -4123.625

$ws.Range("H34").Value = 2273.1462
$ws.Range("I34").Value = 1466.44
$ws.Range("J34").Value = 3533.625
$ws.Range("K34").Value = 1466.44
$ws.Range("L34").Value = 3533.625
$ws.Range("M34").Value = -1264.44
$ws.Range("N34").Value = -3937.625

$ws.Range("H99").Value = 2407.6924
$ws.Range("I99").Value = 1566.6666
$ws.Range("J99").Value = 2517.3914
$ws.Range("K99").Value = 1566.6666
$ws.Range("L99").Value = 2517.3914
$ws.Range("M99").Value = -68.66660000000002
$ws.Range("N99").Value = -5513.3914

$ws.Range("H107").Value = 1277.1154
$ws.Range("I107").Value = 489.64285
$ws.Range("J107").Value = 2195.8333
$ws.Range("K107").Value = 489.64285
$ws.Range("L107").Value = 2195.8333
$ws.Range("M107").Value = 1430.35715
$ws.Range("N107").Value = -6035.8333

$ws.Range("H126").Value = 2407.6924
$ws.Range("I126").Value = 1566.6666
$ws.Range("J126").Value = 2517.3914
$ws.Range("K126").Value = 4699.9998
$ws.Range("L126").Value = 7552.174199999999
$ws.Range("M126").Value = -2229.9998
$ws.Range("N126").Value = -12492.1742

$ws.Range("H132").Value = 3602
$ws.Range("I132").Value = 2450.4614
$ws.Range("J132").Value = 4753.5386
$ws.Range("K132").Value = 7351.3842
$ws.Range("L132").Value = 14260.6158
$ws.Range("M132").Value = -4821.3842
$ws.Range("N132").Value = -19320.6158

$ws.Range("H134").Value = 2762.3044
$ws.Range("I134").Value = 1701.9231
$ws.Range("J134").Value = 4140.8
$ws.Range("K134").Value = 5105.7693
$ws.Range("L134").Value = 12422.4
$ws.Range("M134").Value = -2570.7693
$ws.Range("N134").Value = -17492.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2409
$ws.Range("I132").Value = 1211.4445
$ws.Range("J132").Value = 6001.6665
$ws.Range("K132").Value = 10903.0005
$ws.Range("L132").Value = 54014.9985
$ws.Range("M132").Value = -8373.0005
$ws.Range("N132").Value = -59074.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3870
$ws.Range("I70").Value = 3835
$ws.Range("J70").Value = 3940
$ws.Range("K70").Value = 3835
$ws.Range("L70").Value = 3940
$ws.Range("M70").Value = -3565
$ws.Range("N70").Value = -4480

$ws.Range("H73").Value = 3870
$ws.Range("I73").Value = 3835
$ws.Range("J73").Value = 3940
$ws.Range("K73").Value = 3835
$ws.Range("L73").Value = 3940
$ws.Range("M73").Value = -2899
$ws.Range("N73").Value = -5812

$ws.Range("H97").Value = 2229.8572
$ws.Range("I97").Value = 921.8
$ws.Range("K97").Value = 921.8
$ws.Range("M97").Value = -425.8

$ws.Range("H102").Value = 2651.4443
$ws.Range("I102").Value = 1768.8334
$ws.Range("J102").Value = 4416.6665
$ws.Range("K102").Value = 1768.8334
$ws.Range("L102").Value = 4416.6665
$ws.Range("M102").Value = -146.8334
$ws.Range("N102").Value = -7660.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2171.1765
$ws.Range("I7").Value = 1574.375
$ws.Range("J7").Value = 2701.6667
$ws.Range("K7").Value = 1574.375
$ws.Range("L7").Value = 2701.6667
$ws.Range("M7").Value = -1462.375
$ws.Range("N7").Value = -2925.6667

$ws.Range("H25").Value = 70008
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 70008
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 70008
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -70468

$ws.Range("H40").Value = 2453.818
$ws.Range("I40").Value = 1632
$ws.Range("J40").Value = 3440
$ws.Range("K40").Value = 1632
$ws.Range("L40").Value = 3440
$ws.Range("M40").Value = -1496
$ws.Range("N40").Value = -3712

$ws.Range("H126").Value = 2171.1765
$ws.Range("I126").Value = 1574.375
$ws.Range("J126").Value = 2701.6667
$ws.Range("K126").Value = 4723.125
$ws.Range("L126").Value = 8105.000100000001
$ws.Range("M126").Value = -2253.125
$ws.Range("N126").Value = -13045.0001

$ws.Range("H132").Value = 2912.9666
$ws.Range("I132").Value = 1852.5625
$ws.Range("J132").Value = 4124.857
$ws.Range("K132").Value = 5557.6875
$ws.Range("L132").Value = 12374.571
$ws.Range("M132").Value = -3027.6875
$ws.Range("N132").Value = -17434.571

$ws.Range("H136").Value = 2042.4
$ws.Range("I136").Value = 1583.0588
$ws.Range("J136").Value = 3018.5
$ws.Range("K136").Value = 4749.1764
$ws.Range("L136").Value = 9055.5
$ws.Range("M136").Value = -2199.1764
$ws.Range("N136").Value = -14155.5

Write-Output "applied changes"